$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.Goto($ws.Range("K16"), $true)
